# CSmatchSht added optional parameter:
# the hard-coded DirDBs folder path stored in B1 is simplified from the
# author's full local path down to a short default "C:\DBs\" (consistent
# with that path now being an optional parameter with a sensible default).
# The active selection on Sheet1 also moves from B4 back to B1, the cell
# that holds this path.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("B1").Value = "C:\DBs\"

$ws1.Range("B1").Select()
